$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Step 1: insert a new row at row 4 (shifts data down), set Nodes=6
$ws.Rows.Item(4).Insert()
$ws.Range("A4").Value = 6

# Step 2: resize/extend table to include 2 new columns (H, I) and the extra row
$lo = $ws.ListObjects.Item("Table1")
$lo.Resize($ws.Range("A3:I11"))

# Step 3: set header names via cell values (which drives the table column names)
$ws.Range("H3").Value = "+90 Day + Improved History Data Structure"
$ws.Range("I3").Value = "+HD"

Write-Host "Range: " $lo.Range.Address()

$ws.Range("H4").Value = 232
$ws.Range("G3:G4").Copy()
$ws.Range("H3:I4").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H7").Value = 442
$ws.Range("H10").Value = 663

$ws.Columns.Item(3).ColumnWidth = 20.33203125
$ws.Columns.Item(4).ColumnWidth = 19.5
$ws.Columns.Item(5).ColumnWidth = 21.33203125
$ws.Columns.Item(6).ColumnWidth = 24.33203125
$ws.Columns.Item(7).ColumnWidth = 22
$ws.Columns.Item(8).ColumnWidth = 21.83203125
$ws.Columns.Item(9).ColumnWidth = 22


